$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ground data values (B:F) for rows 2-9
$data = @{
    2 = @(369.2764892578125, 0.277, 0.2587000131607056, 0.7822999954223633, -0.3278999924659729)
    3 = @(438.6117858886719, 0.3509, 0.2982, 1.218999981880188, -0.9613999724388123)
    4 = @(269.6650085449219, 0.315, 0.2751, 1.045600056648254, -0.06830000132322311)
    5 = @(171.4759979248047, 0.1809, 0.1707, 0.6610000133514404, -0.1918999999761581)
    6 = @(-74.42340087890625, -0.0589, -0.0435, 0.4988999962806702, -0.5260000228881836)
    7 = @(-73.51029968261719, -0.074, -0.05979999899864197, 0.3449000120162964, -0.5212000012397766)
    8 = @(380.5320129394531, 0.341, 0.3268, 0.767799973487854, -0.2372999936342239)
    9 = @(1481.627685546875, 0.1918, 0.1691, 1.218999981880188, -0.9613999724388123)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
}
